$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header rows: the reporting window rolled forward by one year/period.
# Column D..H previously covered 1396..1400 (and matching publish dates);
# now it covers 1397..1401, i.e. each column shows what used to be one
# column to its right, plus a brand-new rightmost (H) column.
# ---------------------------------------------------------------------------
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

$ws.Range("D9").Value = "1399-04-16 (8)"
$ws.Range("E9").Value = "1400-04-20 (8)"
$ws.Range("F9").Value = "1401-04-11 (8)"
$ws.Range("G9").Value = "1402-02-28 (7)"
$ws.Range("H9").Value = "1402-02-28"

# ---------------------------------------------------------------------------
# Data rows: same roll-forward for every figure row - the old E:H values
# slide left into D:G, and H gets the newly published figure.
# ---------------------------------------------------------------------------

function Set-Row($r, $d, $e, $f, $g, $h) {
    $ws.Range("D$r").Value = $d
    $ws.Range("E$r").Value = $e
    $ws.Range("F$r").Value = $f
    $ws.Range("G$r").Value = $g
    $ws.Range("H$r").Value = $h
}

Set-Row 12 219096   1105713  529002   1578706  3121725
Set-Row 13 -28934   -48284   -190051  -439855  -191153
Set-Row 14 190162   1057429  338951   1138851  2930572

Set-Row 16 1        4000     9797     0        140
Set-Row 17 -69919   -348533  -663711  -2480531 -1997803
Set-Row 18 0        0        0        0        0
Set-Row 19 0        0        0        0        0
Set-Row 20 -3219    -9427    0        -44118   0
Set-Row 21 0        0        0        0        0
Set-Row 22 0        0        0        0        0
Set-Row 23 0        0        0        120000   0
Set-Row 24 -3954    -5009    -122744  -28432   0
Set-Row 25 0        0        0        0        0
Set-Row 26 0        0        0        0        0
Set-Row 27 0        0        0        0        0
Set-Row 28 0        0        0        0        0
Set-Row 29 0        0        0        0        0
Set-Row 30 258      500      924      5006     1101
Set-Row 31 0        23595    37641    39492    27229
Set-Row 32 -76833   -334874  -738093  -2388583 -1969333
Set-Row 33 113329   722555   -399142  -1249732 961239

Set-Row 35 0        0        0        0        0

# Row 36 keeps the "-" placeholder text in D:F (shifted from E:G), but the
# newly introduced G/H columns report an actual numeric 0 instead of "-".
$ws.Range("D36").Value = "-"
$ws.Range("E36").Value = "-"
$ws.Range("F36").Value = "-"
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 0

Set-Row 37 0        0        331949   266878   0
Set-Row 38 0        0        -681870  -597522  0
Set-Row 39 185849   271776   1033311  2946045  4989955
Set-Row 40 -94318   -80545   -37289   -775638  -3432193
Set-Row 41 -65581   -70011   -49080   -120366  -304813
Set-Row 42 0        0        0        0        0
Set-Row 43 0        0        0        0        0
Set-Row 44 0        0        0        0        0
Set-Row 45 0        0        0        0        0
Set-Row 46 0        0        0        0        0
Set-Row 47 0        0        0        0        0
Set-Row 48 0        0        0        0        0
Set-Row 49 0        0        0        0        0
Set-Row 50 -260     -213717  -547416  -319555  -759822
Set-Row 51 25690    -92497   160661   1399842  493127
Set-Row 52 139019   630058   -238481  150110   1454366
Set-Row 53 54337    192976   823607   585849   737383
Set-Row 54 -380     573      723      1424     3347
Set-Row 55 192976   823607   585849   737383   2195096
Set-Row 56 206539   277721   801577   357836   0
